$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for columns A (Caso) and B (F. De Reclamo) so that
# numeric-looking / date-looking strings are stored as text, matching the source data.
$ws.Range("A19:B32").NumberFormat = "@"

# Row 19
$ws.Range("A19").Value = '7037'
$ws.Range("B19").Value = '8/25/2025'
$ws.Range("C19").Value = 'COCHABAMBA 2556'
$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 'ICD30508273'
$ws.Range("F19").Value = 'Optical Power'
$ws.Range("G19").Value = 'Pendiente'
$ws.Range("H19").Value = 'Tendido a baja altura y cortados'
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = '{"direccionesNormalizadas": [{"altura": 2556, "cod_calle": 3134, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.401027", "y": "-34.625071"}, "direccion": "COCHABAMBA 2556, CABA", "nombre_calle": "COCHABAMBA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K19").Value = -58.401027
$ws.Range("L19").Value = -34.625071
$ws.Range("M19").Value = 'San Telmo'
$ws.Range("N19").Value = 'Capital Sur'

# Row 20
$ws.Range("A20").Value = '3999'
$ws.Range("B20").Value = '8/29/2025'
$ws.Range("C20").Value = 'COCHABAMBA 424'
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 'ICD30532721'
$ws.Range("F20").Value = 'Optical Power'
$ws.Range("G20").Value = 'Pendiente'
$ws.Range("H20").Value = 'Cable en panza'
$ws.Range("I20").Value = 1
$ws.Range("J20").Value = '{"direccionesNormalizadas": [{"altura": 424, "cod_calle": 3134, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.371599", "y": "-34.623098"}, "direccion": "COCHABAMBA 424, CABA", "nombre_calle": "COCHABAMBA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K20").Value = -58.371599
$ws.Range("L20").Value = -34.623098
$ws.Range("M20").Value = 'San Telmo'
$ws.Range("N20").Value = 'Capital Sur'

# Row 21
$ws.Range("A21").Value = '3797'
$ws.Range("B21").Value = '8/29/2025'
$ws.Range("C21").Value = 'CORDOBA AV. 2633'
$ws.Range("D21").Value = 2
$ws.Range("E21").Value = 'ICD30593883'
$ws.Range("F21").Value = 'Optical Power'
$ws.Range("G21").Value = 'Pendiente'
$ws.Range("H21").Value = 'Tendido a baja altura obstaculiza contenedores'
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = '{"direccionesNormalizadas": [{"altura": 2633, "cod_calle": 3165, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.404367", "y": "-34.598010"}, "direccion": "CORDOBA AV. 2633, CABA", "nombre_calle": "CORDOBA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K21").Value = -58.404367
$ws.Range("L21").Value = -34.59801
$ws.Range("M21").Value = 'Almagro'
$ws.Range("N21").Value = 'Capital Sur'

# Row 22
$ws.Range("A22").Value = '3798'
$ws.Range("B22").Value = '8/29/2025'
$ws.Range("C22").Value = 'RIVADAVIA AV. 1559'
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 'ICD30593920'
$ws.Range("F22").Value = 'Optical Power'
$ws.Range("G22").Value = 'Pendiente'
$ws.Range("H22").Value = 'Tendido a baja altura obstaculiza contenedores'
$ws.Range("I22").Value = 1
$ws.Range("J22").Value = '{"direccionesNormalizadas": [{"altura": 1559, "cod_calle": 19046, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.388501", "y": "-34.608971"}, "direccion": "RIVADAVIA AV. 1559, CABA", "nombre_calle": "RIVADAVIA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K22").Value = -58.388501
$ws.Range("L22").Value = -34.608971
$ws.Range("M22").Value = 'San Telmo'
$ws.Range("N22").Value = 'Capital Sur'

# Row 23
$ws.Range("A23").Value = '3878'
$ws.Range("B23").Value = '8/29/2025'
$ws.Range("C23").Value = 'JURAMENTO 3020'
$ws.Range("D23").Value = 13
$ws.Range("E23").Value = 'ICD30574933'
$ws.Range("F23").Value = 'Optical Power'
$ws.Range("G23").Value = 'Pendiente'
$ws.Range("H23").Value = 'Cable en panza'
$ws.Range("I23").Value = 1
$ws.Range("J23").Value = '{"direccionesNormalizadas": [{"altura": 3020, "cod_calle": 10017, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.462535", "y": "-34.565542"}, "direccion": "JURAMENTO 3020, CABA", "nombre_calle": "JURAMENTO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K23").Value = -58.462535
$ws.Range("L23").Value = -34.565542
$ws.Range("M23").Value = 'Colegiales'
$ws.Range("N23").Value = 'Capital Norte'

# Row 24
$ws.Range("A24").Value = '4101'
$ws.Range("B24").Value = '9/2/2025'
$ws.Range("C24").Value = 'SOLER 4197'
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 'ICD30626824'
$ws.Range("F24").Value = 'Optical Power'
$ws.Range("G24").Value = 'Pendiente'
$ws.Range("H24").Value = 'Tendido a baja altura'
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = '{"direccionesNormalizadas": [{"altura": 4197, "cod_calle": 20104, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.420092", "y": "-34.590101"}, "direccion": "SOLER 4197, CABA", "nombre_calle": "SOLER", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K24").Value = -58.420092
$ws.Range("L24").Value = -34.590101
$ws.Range("M24").Value = 'Palermo'
$ws.Range("N24").Value = 'Capital Sur'

# Row 25
$ws.Range("A25").Value = '7160'
$ws.Range("B25").Value = '9/5/2025'
$ws.Range("C25").Value = 'JUNIN 234'
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 'ICD30682500'
$ws.Range("F25").Value = 'Optical Power'
$ws.Range("G25").Value = 'Pendiente'
$ws.Range("H25").Value = 'Cables sueltos y desprolijos'
$ws.Range("I25").Value = 1
$ws.Range("J25").Value = '{"direccionesNormalizadas": [{"altura": 234, "cod_calle": 10015, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.396557", "y": "-34.606731"}, "direccion": "JUNIN 234, CABA", "nombre_calle": "JUNIN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K25").Value = -58.396557
$ws.Range("L25").Value = -34.606731
$ws.Range("M25").Value = 'Almagro'
$ws.Range("N25").Value = 'Capital Sur'

# Row 26
$ws.Range("A26").Value = '7194'
$ws.Range("B26").Value = '9/9/2025'
$ws.Range("C26").Value = 'RIVADAVIA AV. 6381'
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 'Pendiente ADM'
$ws.Range("F26").Value = 'Optical Power'
$ws.Range("G26").Value = 'Pendiente'
$ws.Range("H26").Value = 'Tendido a baja altura'
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = '{"direccionesNormalizadas": [{"altura": 6381, "cod_calle": 19046, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.456364", "y": "-34.626443"}, "direccion": "RIVADAVIA AV. 6381, CABA", "nombre_calle": "RIVADAVIA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K26").Value = -58.456364
$ws.Range("L26").Value = -34.626443
$ws.Range("M26").Value = 'Boedo'
$ws.Range("N26").Value = 'Capital Sur'

# Row 27
$ws.Range("A27").Value = '7198'
$ws.Range("B27").Value = '9/9/2025'
$ws.Range("C27").Value = 'PERU 1305'
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 'Pendiente ADM'
$ws.Range("F27").Value = 'Optical Power'
$ws.Range("G27").Value = 'Pendiente'
$ws.Range("H27").Value = 'Cable en panza'
$ws.Range("I27").Value = 1
$ws.Range("J27").Value = '{"direccionesNormalizadas": [{"altura": 1305, "cod_calle": 17071, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.374174", "y": "-34.623403"}, "direccion": "PERU 1305, CABA", "nombre_calle": "PERU", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K27").Value = -58.374174
$ws.Range("L27").Value = -34.623403
$ws.Range("M27").Value = 'San Telmo'
$ws.Range("N27").Value = 'Capital Sur'

# Row 28
$ws.Range("A28").Value = '7128'
$ws.Range("B28").Value = '9/9/2025'
$ws.Range("C28").Value = 'ARENALES 1928'
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 'Pendiente ADM'
$ws.Range("F28").Value = 'Optical Power'
$ws.Range("G28").Value = 'Pendiente'
$ws.Range("H28").Value = 'Tendido colgando y desordenado se solicita emprolijarlo'
$ws.Range("I28").Value = 1
$ws.Range("J28").Value = '{"direccionesNormalizadas": [{"altura": 1928, "cod_calle": 1104, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.394943", "y": "-34.594725"}, "direccion": "ARENALES 1928, CABA", "nombre_calle": "ARENALES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K28").Value = -58.394943
$ws.Range("L28").Value = -34.594725
$ws.Range("M28").Value = 'Recoleta'
$ws.Range("N28").Value = 'Capital Sur'

# Row 29
$ws.Range("A29").Value = '3740'
$ws.Range("B29").Value = '9/9/2025'
$ws.Range("C29").Value = 'ALCARAZ 4982'
$ws.Range("D29").Value = 10
$ws.Range("E29").Value = 'ICD30722914'
$ws.Range("F29").Value = 'Optical Power'
$ws.Range("G29").Value = 'Pendiente'
$ws.Range("H29").Value = 'Tendido a baja altura'
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = '{"direccionesNormalizadas": [{"altura": 4982, "cod_calle": 1035, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.506496", "y": "-34.621373"}, "direccion": "ALCARAZ 4982, CABA", "nombre_calle": "ALCARAZ", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K29").Value = -58.506496
$ws.Range("L29").Value = -34.621373
$ws.Range("M29").Value = 'Devoto'
$ws.Range("N29").Value = 'Capital Norte'

# Row 30
$ws.Range("A30").Value = '3748'
$ws.Range("B30").Value = '9/9/2025'
$ws.Range("C30").Value = 'JUSTO, JUAN B. AV. 7947'
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = 'ICD30722441'
$ws.Range("F30").Value = 'Optical Power'
$ws.Range("G30").Value = 'Pendiente'
$ws.Range("H30").Value = 'Caja de empalme colgando'
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = '{"direccionesNormalizadas": [{"altura": 7947, "cod_calle": 10018, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.500866", "y": "-34.632009"}, "direccion": "JUSTO, JUAN B. AV. 7947, CABA", "nombre_calle": "JUSTO, JUAN B. AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K30").Value = -58.500866
$ws.Range("L30").Value = -34.632009
$ws.Range("M30").Value = 'Devoto'
$ws.Range("N30").Value = 'Capital Norte'

# Row 31
$ws.Range("A31").Value = '4003'
$ws.Range("B31").Value = '9/9/2025'
$ws.Range("C31").Value = 'CABRERA, JOSE A. 5612'
$ws.Range("D31").Value = 14
$ws.Range("E31").Value = 'Pendiente ADM'
$ws.Range("F31").Value = 'Optical Power'
$ws.Range("G31").Value = 'Pendiente'
$ws.Range("H31").Value = 'Hay un cable colgando ver fotos'
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = '{"direccionesNormalizadas": [{"altura": 5612, "cod_calle": 3009, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.438309", "y": "-34.585646"}, "direccion": "CABRERA, JOSE A. 5612, CABA", "nombre_calle": "CABRERA, JOSE A.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K31").Value = -58.438309
$ws.Range("L31").Value = -34.585646
$ws.Range("M31").Value = 'Palermo'
$ws.Range("N31").Value = 'Capital Sur'

# Row 32
$ws.Range("A32").Value = '6467'
$ws.Range("B32").Value = '9/10/2025'
$ws.Range("C32").Value = 'ASUNCION 2540'
$ws.Range("D32").Value = 15
$ws.Range("E32").Value = 'Pendiente ADM'
$ws.Range("F32").Value = 'Optical Power'
$ws.Range("G32").Value = 'Pendiente'
$ws.Range("H32").Value = 'Cable con baja altura'
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = '{"direccionesNormalizadas": [{"altura": 2540, "cod_calle": 1131, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.492808", "y": "-34.589594"}, "direccion": "ASUNCION 2540, CABA", "nombre_calle": "ASUNCION", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K32").Value = -58.492808
$ws.Range("L32").Value = -34.589594
$ws.Range("M32").Value = 'Paternal'
$ws.Range("N32").Value = 'Capital Norte'
